$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Tagesergebnisse
$ws2 = $wb.Worksheets.Item(2)   # Monatsergebnisse
$ws3 = $wb.Worksheets.Item(3)   # Gesamtergebnis

# ---------------------------------------------------------------------------
# Sheet 1: Tagesergebnisse - insert new "Ein-/Auszahlungen" column at F
# ---------------------------------------------------------------------------
$ws1.Columns.Item(6).Insert() | Out-Null
$ws1.Range("F1").Value = "Ein-/Auszahlungen"
$ws1.Range("F2").Value = 0
$ws1.Range("F2").NumberFormat = "#,##0.00"
$ws1.Columns.Item(6).ColumnWidth = 17

# drop the trailing placeholder rows at the bottom of the sheet
$ws1.Rows.Item(1048574).Delete() | Out-Null
$ws1.Rows.Item(1048574).Delete() | Out-Null
$ws1.Rows.Item(1048574).Delete() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: Monatsergebnisse - insert new "Ein-/Auszahlungen" column at F
# ---------------------------------------------------------------------------
$ws2.Columns.Item(6).Insert() | Out-Null
$ws2.Range("F1").Value = "Ein-/Auszahlungen"
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("F4").Value = 0
$ws2.Range("F5").Value = 0
$ws2.Range("F2:F5").NumberFormat = "#,##0.00"
$ws2.Columns.Item(6).ColumnWidth = 17

# ---------------------------------------------------------------------------
# Sheet 3: Gesamtergebnis - insert new "Ein-/Auszahlungen" column at E
# ---------------------------------------------------------------------------
$ws3.Columns.Item(5).Insert() | Out-Null
$ws3.Range("E1").Value = "Ein-/Auszahlungen"
$ws3.Range("E2").Value = 0
$ws3.Range("E3").Value = 0
$ws3.Range("E2:E3").NumberFormat = "#,##0.00"
$ws3.Columns.Item(5).ColumnWidth = 17

# ---------------------------------------------------------------------------
# Row heights: default row height moves from 15 to 13.8 on every sheet
# ---------------------------------------------------------------------------
$ws1.Rows.Item(1).RowHeight = 13.8
$ws2.Rows.Item(1).RowHeight = 13.8
$ws2.Rows.Item(3).RowHeight = 13.8
$ws2.Rows.Item(4).RowHeight = 13.8
$ws2.Rows.Item(5).RowHeight = 13.8
$ws3.Rows.Item(1).RowHeight = 13.8
$ws3.Rows.Item(2).RowHeight = 13.8

# ---------------------------------------------------------------------------
# Selections + active tab: Tagesergebnisse becomes the active sheet
# ---------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("F11").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("E1").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("F10").Select() | Out-Null

Write-Output "done"
